$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "CreatedAt: 2025-06-17T16:07:43"
$ws.Range("U4").Value = 75
$ws.Range("W4").Value = 31.64
$ws.Range("X4").Value = 29.41
$ws.Range("Y4").Value = 33.95
$ws.Range("Z4").Value = 30.14
$ws.Range("U6").Value = -3.75
$ws.Range("W6").Value = -1.55
$ws.Range("Y6").Value = -0.95
$ws.Range("Z6").Value = -0.51
$ws.Range("X8").Value = 0
$ws.Range("Y8").Value = 0
$ws.Range("T9").Value = 89.37
$ws.Range("U9").Value = 72.84999999999999
$ws.Range("V9").Value = 70.45999999999999
$ws.Range("W9").Value = 31.02
$ws.Range("X9").Value = 29.61
$ws.Range("Y9").Value = 34.45
$ws.Range("Z9").Value = 31.22
$ws.Range("T11").Value = -8.039999999999999
$ws.Range("U11").Value = -5.9
$ws.Range("V11").Value = -5.5
$ws.Range("X11").Value = -0.8
$ws.Range("Y11").Value = -0.45
$ws.Range("Z11").Value = 0.5600000000000001
$ws.Range("X13").Value = 0
$ws.Range("Y13").Value = 0
$ws.Range("T14").Value = 74.58
$ws.Range("U14").Value = 72.84999999999999
$ws.Range("V14").Value = 70.39
$ws.Range("W14").Value = 31.02
$ws.Range("X14").Value = 29.61
$ws.Range("T15").Value = -14.72
$ws.Range("Y15").Value = 156.3
$ws.Range("Z15").Value = 148.39
$ws.Range("T16").Value = -8.130000000000001
$ws.Range("U16").Value = -5.9
$ws.Range("V16").Value = -5.56
$ws.Range("X16").Value = -0.8
$ws.Range("Y16").Value = -0.45
$ws.Range("Z16").Value = 0.59
$ws.Range("X18").Value = 0
$ws.Range("Y18").Value = 0
$ws.Range("T19").Value = 92.59999999999999
$ws.Range("U19").Value = 74.72
$ws.Range("V19").Value = 72.27
$ws.Range("W19").Value = 31.49
$ws.Range("X19").Value = 29.36
$ws.Range("Y19").Value = 34.05
$ws.Range("Z19").Value = 30.41
$ws.Range("T21").Value = -4.82
$ws.Range("U21").Value = -4.03
$ws.Range("V21").Value = -3.69
$ws.Range("W21").Value = -1.7
$ws.Range("Y21").Value = -0.85
$ws.Range("Z21").Value = -0.24
$ws.Range("X23").Value = 0
$ws.Range("Y23").Value = 0
$ws.Range("T24").Value = 92.59999999999999
$ws.Range("U24").Value = 74.72
$ws.Range("V24").Value = 72.27
$ws.Range("W24").Value = 31.49
$ws.Range("X24").Value = 29.36
$ws.Range("Y24").Value = 34.05
$ws.Range("Z24").Value = 30.41
$ws.Range("T26").Value = -4.82
$ws.Range("U26").Value = -4.03
$ws.Range("V26").Value = -3.69
$ws.Range("W26").Value = -1.7
$ws.Range("Y26").Value = -0.85
$ws.Range("Z26").Value = -0.24
$ws.Range("X28").Value = 0
$ws.Range("Y28").Value = 0
$ws.Range("U29").Value = 74.56999999999999
$ws.Range("W29").Value = 31.37
$ws.Range("X29").Value = 29.3
$ws.Range("Y29").Value = 34.15
$ws.Range("Z29").Value = 30.72
$ws.Range("U31").Value = -4.18
$ws.Range("W31").Value = -1.82
$ws.Range("Y31").Value = -0.75
$ws.Range("Z31").Value = 0.06
$ws.Range("X33").Value = 0
$ws.Range("Y33").Value = 0
$ws.Range("U34").Value = 72.12
$ws.Range("V34").Value = 68.92
$ws.Range("W34").Value = 30.5
$ws.Range("X34").Value = 29.73
$ws.Range("Y34").Value = 191.09
$ws.Range("Z34").Value = 180.12
$ws.Range("T35").Value = -14.72
$ws.Range("Y35").Value = 156.3
$ws.Range("Z35").Value = 148.39
$ws.Range("T36").Value = -10.59
$ws.Range("U36").Value = -6.63
$ws.Range("V36").Value = -7.03
$ws.Range("W36").Value = -2.68
$ws.Range("X36").Value = -0.68
$ws.Range("Y36").Value = -0.1
$ws.Range("Z36").Value = 1.08
$ws.Range("X38").Value = 0
$ws.Range("Y38").Value = 0
$ws.Range("U39").Value = 75
$ws.Range("W39").Value = 31.64
$ws.Range("X39").Value = 29.41
$ws.Range("Y39").Value = 33.95
$ws.Range("Z39").Value = 30.14
$ws.Range("U41").Value = -3.75
$ws.Range("W41").Value = -1.55
$ws.Range("Y41").Value = -0.95
$ws.Range("Z41").Value = -0.51
$ws.Range("X43").Value = 0
$ws.Range("Y43").Value = 0
$ws.Range("T44").Value = 95.59999999999999
$ws.Range("U44").Value = 77.51000000000001
$ws.Range("W44").Value = 32.79
$ws.Range("X44").Value = 30.26
$ws.Range("Y44").Value = 34.45
$ws.Range("Z44").Value = 30.26
$ws.Range("T46").Value = -1.82
$ws.Range("U46").Value = -1.24
$ws.Range("W46").Value = -0.39
$ws.Range("X46").Value = -0.15
$ws.Range("Y46").Value = -0.45
$ws.Range("Z46").Value = -0.39
$ws.Range("X48").Value = 0
$ws.Range("Y48").Value = 0
$ws.Range("T49").Value = 92.95999999999999
$ws.Range("U49").Value = 74.29000000000001
$ws.Range("W49").Value = 33.29
$ws.Range("X49").Value = 30.2
$ws.Range("Y49").Value = 34.31
$ws.Range("Z49").Value = 30.44
$ws.Range("T51").Value = -4.46
$ws.Range("U51").Value = -4.46
$ws.Range("X51").Value = -0.21
$ws.Range("Y51").Value = -0.58
$ws.Range("Z51").Value = -0.21
$ws.Range("X53").Value = 0
$ws.Range("Y53").Value = 0
$ws.Range("T54").Value = 90.79000000000001
$ws.Range("U54").Value = 73.73999999999999
$ws.Range("W54").Value = 32.09
$ws.Range("X54").Value = 29.55
$ws.Range("Y54").Value = 33.95
$ws.Range("Z54").Value = 30.11
$ws.Range("T56").Value = -6.63
$ws.Range("U56").Value = -5.01
$ws.Range("W56").Value = -1.09
$ws.Range("X56").Value = -0.86
$ws.Range("Y56").Value = -0.95
$ws.Range("Z56").Value = -0.54
$ws.Range("X58").Value = 0
$ws.Range("Y58").Value = 0
$ws.Range("T59").Value = 99.51000000000001
$ws.Range("U59").Value = 80.52
$ws.Range("V59").Value = 78.14
$ws.Range("W59").Value = 34.14
$ws.Range("X59").Value = 31.32
$ws.Range("Y59").Value = 35.54
$ws.Range("Z59").Value = 31.12
$ws.Range("T61").Value = 2.09
$ws.Range("U61").Value = 1.77
$ws.Range("V61").Value = 2.19
$ws.Range("W61").Value = 0.96
$ws.Range("X61").Value = 0.91
$ws.Range("Y61").Value = 0.64
$ws.Range("Z61").Value = 0.47
$ws.Range("X63").Value = 0
$ws.Range("Y63").Value = 0
$ws.Range("T64").Value = 101.48
$ws.Range("U64").Value = 82.2
$ws.Range("W64").Value = 34.82
$ws.Range("X64").Value = 31.84
$ws.Range("Y64").Value = 36.09
$ws.Range("Z64").Value = 31.57
$ws.Range("T66").Value = 4.06
$ws.Range("U66").Value = 3.45
$ws.Range("W66").Value = 1.64
$ws.Range("X66").Value = 1.43
$ws.Range("Y66").Value = 1.19
$ws.Range("Z66").Value = 0.92
$ws.Range("X68").Value = 0
$ws.Range("Y68").Value = 0
$ws.Range("T69").Value = 102.01
$ws.Range("U69").Value = 82.45999999999999
$ws.Range("W69").Value = 35.08
$ws.Range("X69").Value = 32.18
$ws.Range("Y69").Value = 36.35
$ws.Range("Z69").Value = 31.7
$ws.Range("T71").Value = 4.59
$ws.Range("U71").Value = 3.71
$ws.Range("W71").Value = 1.89
$ws.Range("X71").Value = 1.77
$ws.Range("Y71").Value = 1.45
$ws.Range("Z71").Value = 1.05
$ws.Range("X73").Value = 0
$ws.Range("Y73").Value = 0
$ws.Range("U74").Value = 78.75
$ws.Range("W74").Value = 33.19
$ws.Range("X74").Value = 30.41
$ws.Range("Y74").Value = 34.9
$ws.Range("Z74").Value = 30.65
$ws.Range("X78").Value = 0
$ws.Range("Y78").Value = 0
$ws.Range("U79").Value = 78.75
$ws.Range("W79").Value = 33.19
$ws.Range("X79").Value = 30.41
$ws.Range("Y79").Value = 34.9
$ws.Range("Z79").Value = 30.65
$ws.Range("X83").Value = 0
$ws.Range("Y83").Value = 0
$ws.Range("T84").Value = 88.48
$ws.Range("U84").Value = 72.25
$ws.Range("W84").Value = 32.25
$ws.Range("X84").Value = 29.61
$ws.Range("Y84").Value = 34.01
$ws.Range("Z84").Value = 30.14
$ws.Range("T86").Value = -8.94
$ws.Range("U86").Value = -6.5
$ws.Range("W86").Value = -0.9399999999999999
$ws.Range("X86").Value = -0.8
$ws.Range("Y86").Value = -0.88
$ws.Range("Z86").Value = -0.51
$ws.Range("X88").Value = 0
$ws.Range("Y88").Value = 0
$ws.Range("U89").Value = 74.5
$ws.Range("W89").Value = 31.37
$ws.Range("X89").Value = 29.3
$ws.Range("Y89").Value = 34.15
$ws.Range("Z89").Value = 30.72
$ws.Range("U91").Value = -4.25
$ws.Range("W91").Value = -1.82
$ws.Range("Y91").Value = -0.75
$ws.Range("Z91").Value = 0.06
$ws.Range("X93").Value = 0
$ws.Range("Y93").Value = 0
